# Remove the "is_locked" and "is_enabled" columns (D:E) from the tenant
# import-template header row. These were represented by shared strings
# containing the <%=comment.is_locked_lbl%>... and <%=comment.is_enabled_lbl%>...
# template macros. Deleting the entire columns shifts the following
# columns (order_by, rem) left and removes the now-unused shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1:E1").EntireColumn.Delete()
